$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -1
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = 1
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -4
